$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for data rows 2-10 from 45184 to 45185
$ws.Range("C2:C10").Value = 45185
